$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("control_panel")

# The "country" selector (row 1: country/Fiji) moves to the GUI, so the
# control_panel sheet no longer needs it - drop the row, shifting the
# age_breakpoints row (and the two blank formatted rows below it) up by one.
$ws.Rows.Item(1).Delete()

# The dropdown_lists sheet only existed to back the old in-sheet country
# dropdown; it is no longer referenced, so remove it entirely.
$dl = $wb.Worksheets.Item("dropdown_lists")
$dl.Delete()

$ws.Range("A1").Select()

